# Refresh the crypto price table (commit: "Updated cryptos list on Fri Jun
# 28 14:48:32 UTC 2024 with GitHub Actions").
#
# - Every coin row gets a refreshed Price (D) / Volume(1h) (E) reading.
# - WrappedEther and ShibaInu (rows 16-17) swap rank order.
# - OKB newly enters the list at row 46, pushing Stacks/Maker/ONDO/Cosmos/
#   InjectiveProtocol (previously rows 46-50) down one row each to 47-51;
#   VeChain (previously row 51) drops off the bottom of the table.
#
# Price/Volume values are written with a leading "'" (Excel's own force-text
# entry prefix) so values like "571.74" or "0.0000172" stay text cells,
# matching the sheet's existing text-formatted Price/Volume columns instead
# of being auto-converted to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.905.71"
$ws.Range('E2').Value = "'  -1.81%  "

$ws.Range('D3').Value = "'3.416.43"
$ws.Range('E3').Value = "'  -1.14%  "

$ws.Range('E4').Value = "'  -0.02%  "

$ws.Range('D5').Value = "'571.74"
$ws.Range('E5').Value = "'  -1.61%  "

$ws.Range('D6').Value = "'141.72"
$ws.Range('E6').Value = "'  -4.55%  "

$ws.Range('D7').Value = "'3.417.42"
$ws.Range('E7').Value = "'  -1.16%  "

$ws.Range('E8').Value = "'  +0.04%  "

$ws.Range('D9').Value = "'0.477"
$ws.Range('E9').Value = "'  +0.37%  "

$ws.Range('D10').Value = "'7.55"
$ws.Range('E10').Value = "'  -1.42%  "

$ws.Range('D11').Value = "'0.125"
$ws.Range('E11').Value = "'  +0.43%  "

$ws.Range('D12').Value = "'0.387"
$ws.Range('E12').Value = "'  -0.61%  "

$ws.Range('D13').Value = "'4.003.90"
$ws.Range('E13').Value = "'  -1.04%  "

$ws.Range('D14').Value = "'28.31"
$ws.Range('E14').Value = "'  +2.04%  "

$ws.Range('E15').Value = "'  +0.16%  "

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = "'0.0000172"
$ws.Range('E16').Value = "'  -1.96%  "

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = "'3.420.71"
$ws.Range('E17').Value = "'  -1.28%  "

$ws.Range('D18').Value = "'61.041.50"
$ws.Range('E18').Value = "'  -1.66%  "

$ws.Range('D19').Value = "'6.35"
$ws.Range('E19').Value = "'  +0.99%  "

$ws.Range('D20').Value = "'14.38"
$ws.Range('E20').Value = "'  +1.31%  "

$ws.Range('D21').Value = "'9.35"
$ws.Range('E21').Value = "'  -2.07%  "

$ws.Range('D22').Value = "'392.86"
$ws.Range('E22').Value = "'  +1.26%  "

$ws.Range('E23').Value = "'  -0.17%  "

$ws.Range('D24').Value = "'73.01"
$ws.Range('E24').Value = "'  +0.73%  "

$ws.Range('D25').Value = "'0.994"
$ws.Range('E25').Value = "'  -2.21%  "

$ws.Range('D26').Value = "'0.0000122"
$ws.Range('E26').Value = "'  -2.12%  "

$ws.Range('D27').Value = "'3.576.22"
$ws.Range('E27').Value = "'  -0.56%  "

$ws.Range('E28').Value = "'  -0.12%  "

$ws.Range('D29').Value = "'7.47"
$ws.Range('E29').Value = "'  -4.40%  "

$ws.Range('D30').Value = "'0.998"
$ws.Range('E30').Value = "'  -0.15%  "

$ws.Range('E31').Value = "'  -1.80%  "

$ws.Range('E32').Value = "'  -0.52%  "

$ws.Range('D33').Value = "'1.44"
$ws.Range('E33').Value = "'  -7.93%  "

$ws.Range('D35').Value = "'23.82"
$ws.Range('E35').Value = "'  -1.03%  "

$ws.Range('E36').Value = "'  -0.74%  "

$ws.Range('D37').Value = "'3.445.34"
$ws.Range('E37').Value = "'  -0.93%  "

$ws.Range('E38').Value = "'  -3.43%  "

$ws.Range('E39').Value = "'  -1.77%  "

$ws.Range('D40').Value = "'167.16"
$ws.Range('E40').Value = "'  +0.43%  "

$ws.Range('D41').Value = "'0.0785"
$ws.Range('E41').Value = "'  -1.26%  "

$ws.Range('D42').Value = "'27.24"
$ws.Range('E42').Value = "'  +3.91%  "

$ws.Range('D43').Value = "'0.797"
$ws.Range('E43').Value = "'  +0.42%  "

$ws.Range('E44').Value = "'  +0.02%  "

$ws.Range('E45').Value = "'  +0.70%  "

$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = "'41.97"
$ws.Range('E46').Value = "'  -0.48%  "

$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = "'1.71"
$ws.Range('E47').Value = "'  -1.50%  "

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = "'2.608.85"
$ws.Range('E48').Value = "'  -0.84%  "

$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = "'1.14"
$ws.Range('E49').Value = "'  -4.54%  "

$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = "'6.96"
$ws.Range('E50').Value = "'  +1.10%  "

$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = "'22.91"
$ws.Range('E51').Value = "'  -4.53%  "
